# "all mlcc are done" - fix up the "Vérifier C73 C74 ..." note (flag the
# word "coté" as a spell-check hit, matching "antialiasing" right next to
# it) and add a fresh "C201 ??? " note below it, same way the existing
# "C201 placé ..." note above it looks.

$d = $word.ActiveDocument

# Namespace-qualified WordprocessingML body fragment, wrapped in the
# single-file "pkg:package" envelope that Range.InsertXML expects -
# InsertXML REPLACES whatever the target range currently spans, so it is
# the right tool both for "rewrite this paragraph's runs" and for
# "splice extra paragraphs in at this collapsed point".
function New-WordPackageXml {
    param([string]$BodyXml)

    return "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">" +
        "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">" +
        "<pkg:xmlData><w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">" +
        "<w:body>$BodyXml</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

# Locate the "Vérifier C73 C74 ..." paragraph by its text rather than a
# hard-coded index.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Vérifier C73 C74*") {
        $targetPara = $p
    }
}

# Rewrite its runs so "coté" is wrapped in spellStart/spellEnd proofErr
# marks, same as "antialiasing" already was.
$fixedRunsXml = (
    '<w:r><w:t xml:space="preserve">Vérifier C73 C74 et miroir du </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>coté</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> des </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>antialiasing</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
)
$targetPara.Range.InsertXML((New-WordPackageXml "<w:p>$fixedRunsXml</w:p>"))

# Re-fetch the (now rewritten) paragraph and insert the new blank /
# "C201 ??? " / blank paragraphs right after it - before the collapsed
# insertion point, so the paragraph that already follows is left
# completely untouched.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Vérifier C73 C74*") {
        $targetPara = $p
    }
}
$afterPara = $targetPara.Next()
$insertAt = $afterPara.Range.Start
$insertionPoint = $d.Range($insertAt, $insertAt)

$newParasXml = '<w:p/><w:p><w:r><w:t xml:space="preserve">C201 ??? </w:t></w:r></w:p><w:p/>'
$insertionPoint.InsertXML((New-WordPackageXml $newParasXml))
